$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 232, shifting existing rows 232:321 down to 233:322
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new record
$ws.Cells.Item(232, 1).Value = 4
$ws.Cells.Item(232, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(232, 3).Value = "Los Lagos"
$ws.Cells.Item(232, 4).Value = 44795
$ws.Cells.Item(232, 5).Value = 10
$ws.Cells.Item(232, 6).Value = 100112037
$ws.Cells.Item(232, 7).Value = "Cebollín"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 70
$ws.Cells.Item(232, 11).Value = 9500
$ws.Cells.Item(232, 12).Value = 10000
$ws.Cells.Item(232, 13).Value = 9750
$ws.Cells.Item(232, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(232, 15).Value = "Región Metropolitana"
$ws.Cells.Item(232, 16).Value = 271
$ws.Cells.Item(232, 17).Value = 36
$ws.Cells.Item(232, 18).Value = "Hortaliza"
